# "Updated test report listener" - refresh the test-run timestamps and a
# couple of comment strings that the listener now emits, and let column E
# (Comment) re-fit the new, longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Execution Time column (D) -------------------------------------------
# Rows 2-14 all share the earlier run's timestamp -> new run timestamp.
$ws.Range("D2:D14").Value = "03/06/2025 01:52:49 PM"
# Row 15/16 (Navigate to Profile Page / Verify profile page loads) share a
# slightly later timestamp from the same run.
$ws.Range("D15:D16").Value = "03/06/2025 01:52:50 PM"
# Row 17 (Logged out successfully) timestamp.
$ws.Range("D17").Value = "03/06/2025 01:52:51 PM"

# --- Comment column (E) ----------------------------------------------------
# The generic "Test was skipped." comment is replaced with a more
# descriptive message for the skipped login tests (rows 3-13).
$ws.Range("E3:E13").Value = "Skipping invalid login test as per config"

# Column E needs to grow to fit the longer comment text.
$ws.Columns.Item(5).ColumnWidth = 34.75
